$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2515.7
$ws.Range("I2").Value = 2739.6667
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 2739.6667
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -2626.6667
$ws.Range("N2").Value = -726
$ws.Range("H9").Value = 398.5
$ws.Range("I9").Value = 355.42856
$ws.Range("J9").Value = 700
$ws.Range("K9").Value = 355.42856
$ws.Range("L9").Value = 700
$ws.Range("M9").Value = -186.42856
$ws.Range("H12").Value = 299.22223
$ws.Range("I12").Value = 274.0625
$ws.Range("J12").Value = 500.5
$ws.Range("K12").Value = 274.0625
$ws.Range("L12").Value = 500.5
$ws.Range("M12").Value = -104.0625
$ws.Range("H17").Value = 1755.5555
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1755.5555
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5266.666499999999
$ws.Range("N17").Value = -5602.666499999999
$ws.Range("H18").Value = 1137
$ws.Range("I18").Value = 681.6667
$ws.Range("J18").Value = 1820
$ws.Range("K18").Value = 681.6667
$ws.Range("L18").Value = 1820
$ws.Range("M18").Value = -397.6667
$ws.Range("N18").Value = -2388
$ws.Range("H32").Value = 2626
$ws.Range("I32").Value = 1918
$ws.Range("J32").Value = 4750
$ws.Range("K32").Value = 1918
$ws.Range("L32").Value = 4750
$ws.Range("M32").Value = -1592
$ws.Range("H34").Value = 11498.25
$ws.Range("I34").Value = 11498.25
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 11498.25
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -11295.25
$ws.Range("H36").Value = 11498.25
$ws.Range("I36").Value = 11498.25
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 11498.25
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -10783.25
$ws.Range("H92").Value = 215.73685
$ws.Range("I92").Value = 217.44444
$ws.Range("J92").Value = 185
$ws.Range("K92").Value = 217.44444
$ws.Range("L92").Value = 185
$ws.Range("M92").Value = 1030.55556
$ws.Range("N92").Value = -2681
$ws.Range("H132").Value = 1890.1666
$ws.Range("I132").Value = 1713.04
$ws.Range("J132").Value = 2775.8
$ws.Range("K132").Value = 5139.12
$ws.Range("L132").Value = 8327.400000000001
$ws.Range("M132").Value = -2609.12
$ws.Range("H135").Value = 2877
$ws.Range("I135").Value = 2877
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 25893
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -23358
$ws.Range("H137").Value = 5113.1924
$ws.Range("I137").Value = 3013.85
$ws.Range("J137").Value = 12111
$ws.Range("K137").Value = 9041.549999999999
$ws.Range("L137").Value = 36333
$ws.Range("M137").Value = -6491.549999999999
$ws.Range("N137").Value = -41433
$ws.Range("H138").Value = 2495.3655
$ws.Range("I138").Value = 1396.6471
$ws.Range("J138").Value = 3029.0286
$ws.Range("K138").Value = 4189.9413
$ws.Range("L138").Value = 9087.085800000001
$ws.Range("M138").Value = 950.0587000000005
$ws.Range("N138").Value = -19367.0858

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15661293
$ws.Range("I32").Value = 17892476
$ws.Range("J32").Value = 43006.25
$ws.Range("K32").Value = 17892476
$ws.Range("L32").Value = 43006.25
$ws.Range("M32").Value = -17892189
$ws.Range("N32").Value = -43580.25
$ws.Range("H45").Value = 29414452
$ws.Range("I45").Value = 50001740
$ws.Range("J45").Value = 4040
$ws.Range("K45").Value = 50001740
$ws.Range("L45").Value = 4040
$ws.Range("M45").Value = -50001363
$ws.Range("H54").Value = 37000
$ws.Range("I54").Value = 39000
$ws.Range("J54").Value = 35000
$ws.Range("K54").Value = 39000
$ws.Range("L54").Value = 35000
$ws.Range("M54").Value = -38231
$ws.Range("N54").Value = -36538
$ws.Range("H61").Value = 50010840
$ws.Range("I61").Value = 45464496
$ws.Range("J61").Value = 62513290
$ws.Range("K61").Value = 45464496
$ws.Range("L61").Value = 62513290
$ws.Range("M61").Value = -45464284
$ws.Range("H63").Value = 5650.357
$ws.Range("I63").Value = 2516.6667
$ws.Range("J63").Value = 8000.625
$ws.Range("K63").Value = 2516.6667
$ws.Range("L63").Value = 8000.625
$ws.Range("M63").Value = -1830.6667
$ws.Range("H66").Value = 5650.357
$ws.Range("I66").Value = 2516.6667
$ws.Range("J66").Value = 8000.625
$ws.Range("K66").Value = 12583.3335
$ws.Range("L66").Value = 40003.125
$ws.Range("M66").Value = -9151.333500000001
$ws.Range("H74").Value = 9290453
$ws.Range("I74").Value = 15628215
$ws.Range("J74").Value = 840104.25
$ws.Range("K74").Value = 15628215
$ws.Range("L74").Value = 840104.25
$ws.Range("M74").Value = -15627341
$ws.Range("H77").Value = 9290453
$ws.Range("I77").Value = 15628215
$ws.Range("J77").Value = 840104.25
$ws.Range("K77").Value = 78141075
$ws.Range("L77").Value = 4200521.25
$ws.Range("M77").Value = -78136707
$ws.Range("H132").Value = 2704.5632
$ws.Range("I132").Value = 1920.1846
$ws.Range("J132").Value = 5022.0454
$ws.Range("K132").Value = 5760.5538
$ws.Range("L132").Value = 15066.1362
$ws.Range("M132").Value = -3230.5538
$ws.Range("N132").Value = -20126.1362
$ws.Range("H136").Value = 50010840
$ws.Range("I136").Value = 45464496
$ws.Range("J136").Value = 62513290
$ws.Range("K136").Value = 136393488
$ws.Range("L136").Value = 187539870
$ws.Range("M136").Value = -136390938

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 271744.47
$ws.Range("I134").Value = 1525.1818
$ws.Range("J134").Value = 2501053.5
$ws.Range("K134").Value = 4575.5454
$ws.Range("L134").Value = 7503160.5
$ws.Range("M134").Value = -2040.5454

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 80000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 80000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 80000
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -81132
$ws.Range("H58").Value = 1789.4147
$ws.Range("I58").Value = 1104.9412
$ws.Range("J58").Value = 5114
$ws.Range("K58").Value = 1104.9412
$ws.Range("L58").Value = 5114
$ws.Range("M58").Value = -901.9412
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1483.8823
$ws.Range("I132").Value = 1634
$ws.Range("J132").Value = 783.3333
$ws.Range("K132").Value = 4902
$ws.Range("L132").Value = 2349.9999
$ws.Range("M132").Value = -2372
$ws.Range("N132").Value = -7409.9999
$ws.Range("H136").Value = 1789.4147
$ws.Range("I136").Value = 1104.9412
$ws.Range("J136").Value = 5114
$ws.Range("K136").Value = 3314.8236
$ws.Range("L136").Value = 15342
$ws.Range("M136").Value = -764.8235999999997

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 861.5
$ws.Range("I11").Value = 34
$ws.Range("J11").Value = 4999
$ws.Range("K11").Value = 102
$ws.Range("L11").Value = 14997
$ws.Range("M11").Value = 38
$ws.Range("N11").Value = -15277
$ws.Range("H26").Value = 199
$ws.Range("I26").Value = 226.66667
$ws.Range("J26").Value = 171.33333
$ws.Range("K26").Value = 680.00001
$ws.Range("L26").Value = 513.99999
$ws.Range("M26").Value = -392.00001
$ws.Range("N26").Value = -1089.99999
$ws.Range("H61").Value = 79.5
$ws.Range("I61").Value = 79.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 238.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -23.5
$ws.Range("N61").ClearContents()
$ws.Range("H93").Value = 13555
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 13555
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 40665
$ws.Range("N93").Value = -44409
$ws.Range("H131").Value = 4820.5737
$ws.Range("I131").Value = 4488.769
$ws.Range("J131").Value = 4910.4375
$ws.Range("K131").Value = 13466.307
$ws.Range("L131").Value = 14731.3125
$ws.Range("M131").Value = -8426.307000000001
$ws.Range("N131").Value = -24811.3125
$ws.Range("H136").Value = 9795.333000000001
$ws.Range("I136").Value = 9795.333000000001
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 29385.999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -24285.999
$ws.Range("H137").Value = 7580
$ws.Range("I137").Value = 8000
$ws.Range("J137").Value = 7486.6665
$ws.Range("K137").Value = 24000
$ws.Range("L137").Value = 22459.9995
$ws.Range("M137").Value = -18900
$ws.Range("N137").Value = -32659.9995
$ws.Range("H138").Value = 4916.6665
$ws.Range("I138").Value = 4500
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 13500
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = -8360
$ws.Range("N138").Value = -25280
$ws.Range("H140").Value = 276979.38
$ws.Range("H141").Value = 441428.56
$ws.Range("I141").Value = 3000000
$ws.Range("J141").Value = 15000
$ws.Range("K141").Value = 9000000
$ws.Range("L141").Value = 45000
$ws.Range("M141").Value = -8994820

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6205.467
$ws.Range("I70").Value = 5037.8
$ws.Range("J70").Value = 6789.3
$ws.Range("K70").Value = 5037.8
$ws.Range("L70").Value = 6789.3
$ws.Range("M70").Value = -4767.8
$ws.Range("H73").Value = 6205.467
$ws.Range("I73").Value = 5037.8
$ws.Range("J73").Value = 6789.3
$ws.Range("K73").Value = 5037.8
$ws.Range("L73").Value = 6789.3
$ws.Range("M73").Value = -4101.8
$ws.Range("H132").Value = 28578466
$ws.Range("I132").Value = 35718492
$ws.Range("J132").Value = 18365.572
$ws.Range("K132").Value = 107155476
$ws.Range("L132").Value = 55096.716
$ws.Range("M132").Value = -107152946
$ws.Range("N132").Value = -60156.716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 75000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 75000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 75000
$ws.Range("N99").Value = -80990
$ws.Range("H106").Value = 37319.4
$ws.Range("I106").Value = 59000
$ws.Range("J106").Value = 31899.25
$ws.Range("K106").Value = 59000
$ws.Range("L106").Value = 31899.25
$ws.Range("M106").Value = -57738
$ws.Range("N106").Value = -34423.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 17666.666
$ws.Range("I15").Value = 17666.666
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 17666.666
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -17378.666
$ws.Range("H20").Value = 29988
$ws.Range("I20").Value = 29988
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 29988
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -29748
$ws.Range("N20").ClearContents()
$ws.Range("H43").Value = 78000
$ws.Range("I43").Value = 76000
$ws.Range("J43").Value = 80000
$ws.Range("K43").Value = 76000
$ws.Range("L43").Value = 80000
$ws.Range("M43").Value = -75851
$ws.Range("N43").Value = -80298
$ws.Range("H49").Value = 33495
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 33495
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 33495
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -33955
$ws.Range("H95").Value = 98999.60000000001
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 98999.60000000001
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 98999.60000000001
$ws.Range("N95").Value = -104491.6
$ws.Range("H119").Value = 44228.25
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 44228.25
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 44228.25
$ws.Range("N119").Value = -53904.25
$ws.Range("H132").Value = 1646.0975
$ws.Range("I132").Value = 1260.3334
$ws.Range("J132").Value = 2698.182
$ws.Range("K132").Value = 3781.0002
$ws.Range("L132").Value = 8094.545999999999
$ws.Range("M132").Value = -1251.0002
